$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The final two columns (D and E) are removed entirely; remaining cells
# shift left automatically.
$ws.Columns.Item(4).Delete()
$ws.Columns.Item(4).Delete()

# Update the remaining header labels.
$ws.Range("B1").Value = "CA_LF"
$ws.Range("C1").Value = "LF_FFR"

# Update the "params" row values.
$ws.Range("B2").Value = 0.8141847678834518
$ws.Range("C2").Value = 0.627195551619356

# Update the "pvalue" row values.
$ws.Range("B3").Value = [double]"3.086993074852273E-05"
$ws.Range("C3").Value = 0

Write-Output "applied edits"
